$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new book entry on row 5 (previously a blank placeholder row).
$ws.Range("A5").Value = "f76eb274"
$ws.Range("B5").Value = "2024-07-25"
$ws.Range("C5").Value = "날아라, 씨앗 폭탄!"
$ws.Range("D5").Value = "이묘신"
$ws.Range("E5").Value = "가로세로그림책"
$ws.Range("F5").Value = 17
$ws.Range("G5").Value = 44
$ws.Range("H5").Value = "한국어"
$ws.Range("I5").Value = "초록개구리"
$ws.Range("J5").Value = "2024-05-25"
$ws.Range("K5").Value = "9791157822973"
$ws.Range("L5").Value = "달달숲 마을엔 나무가 없다. 사람들이 모조리 베어 가서 그루터기만 가득하다. 어느 날, 어른 동물들은 ‘폭탄’을 만들기로 한다. 그 소식을 엿들은 아기 여우는 헐레벌떡 친구들에게 달려간다. 아기 동물들은 어른들이 전쟁을 벌일 거라는 생각에 그 ‘폭탄’을 찾아서 꼭꼭 숨긴다. 그런데… 참 이상해하다. 며칠 뒤 ‘폭탄’을 숨긴 곳에 파릇파릇한 싹이 돋아난 것이다. 이 수상한 폭탄의 정체는 무엇일까?"

# Move the active selection to B2, as recorded in the saved view state.
$ws.Range("B2").Select()
